$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")
$ws.Activate()

# Core model edit: use the exact fraction 2/3 instead of the rounded 0.66
# for the "reduced model" ion-energy fraction, and drop the now-unused
# helper cell next to it.
$ws.Range("D24").Formula = "=2/3*Einj"
$ws.Range("E24").Clear()

# E25 picks up the plain-integer number format that D24/E24 used to carry.
$ws.Range("E25").NumberFormat = "0"

# View-state touch-ups to mirror the re-saved workbook (best effort - the
# exact zoom bookkeeping / auto row-height metrics are cosmetic artifacts
# of the authoring app and aren't reproduced bit-for-bit here).
$excel.ActiveWindow.ScrollRow = 1
$excel.ActiveWindow.ScrollColumn = 1
$excel.ActiveWindow.Zoom = 100
$ws.Range("D34").Select() | Out-Null
